{"js": "// Update the multiplication expressions in the table (old -> new), matching\n// the exact text of each cell's w:t run. Every source expression in this\n// document is unique, so a direct search+replace per pair is unambiguous.\nconst replacements = [\n  [\"43\u00d763=\", \"70\u00d755=\"],\n  [\"14\u00d732=\", \"87\u00d740=\"],\n  [\"54\u00d721=\", \"16\u00d714=\"],\n  [\"50\u00d748=\", \"96\u00d722=\"],\n  [\"62\u00d769=\", \"64\u00d759=\"],\n  [\"12\u00d782=\", \"21\u00d740=\"],\n  [\"73\u00d756=\", \"74\u00d752=\"],\n  [\"46\u00d753=\", \"41\u00d787=\"],\n  [\"50\u00d793=\", \"15\u00d769=\"],\n  [\"34\u00d716=\", \"80\u00d741=\"],\n  [\"42\u00d751=\", \"19\u00d791=\"],\n  [\"12\u00d717=\", \"53\u00d712=\"],\n  [\"80\u00d753=\", \"53\u00d766=\"],\n  [\"78\u00d769=\", \"70\u00d755=\"],\n  [\"45\u00d767=\", \"90\u00d752=\"],\n  [\"20\u00d735=\", \"96\u00d784=\"],\n  [\"29\u00d748=\", \"36\u00d773=\"],\n  [\"41\u00d746=\", \"40\u00d787=\"],\n  [\"75\u00d713=\", \"18\u00d733=\"],\n  [\"97\u00d756=\", \"13\u00d775=\"],\n  [\"69\u00d782=\", \"93\u00d759=\"],\n  [\"33\u00d729=\", \"90\u00d759=\"],\n  [\"40\u00d772=\", \"31\u00d787=\"],\n  [\"95\u00d730=\", \"81\u00d713=\"],\n  [\"58\u00d730=\", \"61\u00d776=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Update the multiplication expressions in the table (old -> new). Every\n# source expression in this document is unique, so a single Find/Replace\n# (ReplaceOne) per pair, run in document order, is unambiguous.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"43\u00d763=\", \"70\u00d755=\"),\n    @(\"14\u00d732=\", \"87\u00d740=\"),\n    @(\"54\u00d721=\", \"16\u00d714=\"),\n    @(\"50\u00d748=\", \"96\u00d722=\"),\n    @(\"62\u00d769=\", \"64\u00d759=\"),\n    @(\"12\u00d782=\", \"21\u00d740=\"),\n    @(\"73\u00d756=\", \"74\u00d752=\"),\n    @(\"46\u00d753=\", \"41\u00d787=\"),\n    @(\"50\u00d793=\", \"15\u00d769=\"),\n    @(\"34\u00d716=\", \"80\u00d741=\"),\n    @(\"42\u00d751=\", \"19\u00d791=\"),\n    @(\"12\u00d717=\", \"53\u00d712=\"),\n    @(\"80\u00d753=\", \"53\u00d766=\"),\n    @(\"78\u00d769=\", \"70\u00d755=\"),\n    @(\"45\u00d767=\", \"90\u00d752=\"),\n    @(\"20\u00d735=\", \"96\u00d784=\"),\n    @(\"29\u00d748=\", \"36\u00d773=\"),\n    @(\"41\u00d746=\", \"40\u00d787=\"),\n    @(\"75\u00d713=\", \"18\u00d733=\"),\n    @(\"97\u00d756=\", \"13\u00d775=\"),\n    @(\"69\u00d782=\", \"93\u00d759=\"),\n    @(\"33\u00d729=\", \"90\u00d759=\"),\n    @(\"40\u00d772=\", \"31\u00d787=\"),\n    @(\"95\u00d730=\", \"81\u00d713=\"),\n    @(\"58\u00d730=\", \"61\u00d776=\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $newText\n    $find.Forward = $true\n    $find.Wrap = 0\n    $find.Format = $false\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n\n    $find.Execute($oldText, $false, $true, $false, $false, $false, $true, 1, $false, $newText, 1)\n}\n"}
